# Stock App Project Planning Notes - Todo list update
# - Mark several items done by prefixing "X "
# - Remove two stale/duplicate todo items
# - Bold the "Add ability to specify date range for chart data" item
# - Insert two new todo items ("About Page", "Logout")

$d = $word.ActiveDocument

function Get-ParaByText($text) {
    $rng = $d.Content
    $rng.Find.ClearFormatting()
    $found = $rng.Find.Execute($text, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Text not found: $text"
    }
    return $rng.Paragraphs(1)
}

function Add-XPrefix($text) {
    $p = Get-ParaByText($text)
    $r = $p.Range.Duplicate()
    $r.Collapse(1)
    $r.InsertBefore("X ")
}

# 1) "API fetch to get current stock price ..." -> prefix "X "
Add-XPrefix("API fetch to get current stock price (use the Stock Price: Quote API)")

# 2) "Add to watchlist view?" -> prefix "X "
Add-XPrefix("Add to watchlist view?")

# 3) Remove stale items: "Finish watchlist detail" and
#    "Add function to show current stock price (future)"
$p = Get-ParaByText("Finish watchlist detail")
$p.Range.Delete()

$p = Get-ParaByText("Add function to show current stock price (future)")
$p.Range.Delete()

# 4) "Add to watchlist button in company view?" -> prefix "X "
Add-XPrefix("Add to watchlist button in company view?")

# 5) "Format company data (numbers are raw currently)" -> prefix "X "
Add-XPrefix("Format company data (numbers are raw currently)")

# 6) "Load latest stock info in company view" -> prefix "X "
Add-XPrefix("Load latest stock info in company view")

# 7) Bold "Add ability to specify date range for chart data"
$p = Get-ParaByText("Add ability to specify date range for chart data")
$p.Range.Font.Bold = 1
$p.Range.Font.BoldBi = 1

# 8) Insert "About Page" before "Homepage"
$p = Get-ParaByText("Homepage")
$p.Range.InsertParagraphBefore()
$newPara = Get-ParaByText("Homepage").Previous()
$newPara.Range.Text = "About Page"

# 9) Insert "Logout" after "New user form / route (/signup)"
$p = Get-ParaByText("New user form / route (/signup)")
$p.Range.InsertParagraphAfter()
$newPara = Get-ParaByText("New user form / route (/signup)").Next()
$newPara.Range.Text = "Logout"

Write-Output "Done"
